# Update the error-table statistics (columns B:G, rows 2:11) with the
# refreshed ifoCAST GVA component-analysis preprocessing results.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = -0.00515214134517856
$ws.Range("C2").Value = 0.6651924495933969
$ws.Range("D2").Value = 0.7443617204222924
$ws.Range("E2").Value = 0.8627640004209103
$ws.Range("F2").Value = 0.8863899424499349
$ws.Range("G2").Value = 19

$ws.Range("B3").Value = -0.04706953193124818
$ws.Range("C3").Value = 0.7263783780799662
$ws.Range("D3").Value = 0.9076748761767716
$ws.Range("E3").Value = 0.9527197259303344
$ws.Range("F3").Value = 0.9791433241278349
$ws.Range("G3").Value = 18

$ws.Range("B4").Value = 0.008980548322284209
$ws.Range("C4").Value = 0.7068694245681296
$ws.Range("D4").Value = 0.7853984072668317
$ws.Range("E4").Value = 0.8862270630413133
$ws.Range("F4").Value = 0.9134550436843235
$ws.Range("G4").Value = 17

$ws.Range("B5").Value = 0.130843318123659
$ws.Range("C5").Value = 0.6351380739751387
$ws.Range("D5").Value = 0.6444494568402878
$ws.Range("E5").Value = 0.8027760938395511
$ws.Range("F5").Value = 0.8180167776226784
$ws.Range("G5").Value = 16

$ws.Range("B6").Value = 0.07882986512450849
$ws.Range("C6").Value = 0.6136611455857488
$ws.Range("D6").Value = 0.5739360866201003
$ws.Range("E6").Value = 0.7575856958919567
$ws.Range("F6").Value = 0.779918909922613
$ws.Range("G6").Value = 15

$ws.Range("B7").Value = 0.1364943803647256
$ws.Range("C7").Value = 0.5668254508253494
$ws.Range("D7").Value = 0.6024723940093601
$ws.Range("E7").Value = 0.7761909520275021
$ws.Range("F7").Value = 0.7929392009836149
$ws.Range("G7").Value = 14

$ws.Range("B8").Value = 0.1194123702819532
$ws.Range("C8").Value = 0.6296987908481204
$ws.Range("D8").Value = 0.6261569364748241
$ws.Range("E8").Value = 0.7913007876116541
$ws.Range("F8").Value = 0.8141800113140473
$ws.Range("G8").Value = 13

$ws.Range("B9").Value = 0.2385215453226533
$ws.Range("C9").Value = 0.5677170375037683
$ws.Range("D9").Value = 0.505199978683133
$ws.Range("E9").Value = 0.7107742107611481
$ws.Range("F9").Value = 0.6993301608877609
$ws.Range("G9").Value = 12

$ws.Range("B10").Value = 0.1002722471476047
$ws.Range("C10").Value = 0.3965624711381974
$ws.Range("D10").Value = 0.2295019047709758
$ws.Range("E10").Value = 0.4790635706991044
$ws.Range("F10").Value = 0.49131671999357
$ws.Range("G10").Value = 11

$ws.Range("B11").Value = 0.1714703583411825
$ws.Range("C11").Value = 0.390033206010506
$ws.Range("D11").Value = 0.2698760920897213
$ws.Range("E11").Value = 0.519495998146012
$ws.Range("F11").Value = 0.5169074796862885
$ws.Range("G11").Value = 10
